$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price / percentage columns (D, E) must stay as literal
# text, matching the workbook's inlineStr cell storage. Force text format
# before assigning so Excel does not auto-convert "331.44" / "0.46%" into
# real numbers.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "E18",
    "E19",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "D46",
    "E46",
    "E47",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "331.44"
$ws.Range("E2").Value = "0.46%"
$ws.Range("D3").Value = "45.44"
$ws.Range("E3").Value = "2.92%"
$ws.Range("D4").Value = "5.611"
$ws.Range("E4").Value = "2.07%"
$ws.Range("D5").Value = "0.08349"
$ws.Range("E5").Value = "4.35%"
$ws.Range("D6").Value = "2.088"
$ws.Range("E6").Value = "5.49%"
$ws.Range("D7").Value = "0.9657"
$ws.Range("E7").Value = "1.51%"
$ws.Range("E8").Value = "-0.85%"
$ws.Range("D9").Value = "0.1173"
$ws.Range("E9").Value = "5.69%"
$ws.Range("D10").Value = "0.1918"
$ws.Range("E10").Value = "0.43%"
$ws.Range("E11").Value = "7.42%"
$ws.Range("D12").Value = "0.09856"
$ws.Range("E12").Value = "-1.24%"
$ws.Range("D13").Value = "0.04614"
$ws.Range("E13").Value = "-3.41%"
$ws.Range("E14").Value = "-0.28%"
$ws.Range("D15").Value = "0.001278"
$ws.Range("E15").Value = "0.56%"
$ws.Range("D16").Value = "0.006077"
$ws.Range("E16").Value = "2.38%"
$ws.Range("D17").Value = "3.377"
$ws.Range("E17").Value = "0.21%"
$ws.Range("E18").Value = "1.13%"
$ws.Range("E19").Value = "-4.25%"
$ws.Range("D21").Value = "0.2879"
$ws.Range("E21").Value = "11.40%"
$ws.Range("D22").Value = "0.04184"
$ws.Range("E22").Value = "2.60%"
$ws.Range("D23").Value = "0.001316"
$ws.Range("E23").Value = "3.45%"
$ws.Range("D24").Value = "0.004552"
$ws.Range("E24").Value = "4.22%"
$ws.Range("D25").Value = "0.0001303"
$ws.Range("E25").Value = "8.62%"
$ws.Range("D26").Value = "0.0003749"
$ws.Range("E26").Value = "0.16%"
$ws.Range("D38").Value = "0.02716"
$ws.Range("E38").Value = "4.63%"
$ws.Range("D39").Value = "0.05757"
$ws.Range("E39").Value = "0.44%"
$ws.Range("D40").Value = "0.007867"
$ws.Range("E40").Value = "4.26%"
$ws.Range("D41").Value = "0.1434"
$ws.Range("E41").Value = "2.37%"
$ws.Range("D42").Value = "0.007269"
$ws.Range("E42").Value = "-1.24%"
$ws.Range("D43").Value = "0.002025"
$ws.Range("E43").Value = "0.51%"
$ws.Range("D44").Value = "0.009091"
$ws.Range("E44").Value = "9.03%"
$ws.Range("D45").Value = "0.3546"
$ws.Range("D46").Value = "0.00007115"
$ws.Range("E46").Value = "-0.11%"
$ws.Range("E47").Value = "0.27%"
$ws.Range("E48").Value = "0.31%"
$ws.Range("D49").Value = "0.003526"
$ws.Range("E49").Value = "-0.85%"
$ws.Range("D50").Value = "0.003506"
$ws.Range("E50").Value = "-0.64%"
$ws.Range("D51").Value = "0.00002105"
$ws.Range("E51").Value = "0.27%"

# Coin name / link columns (B, C) are plain text already.
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("B50").Value = "CoinbaseStockToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
